# Apply the "F-Measure" section addition to Result_Values.xlsx
# - Update the existing Precision/ROS-AUC style numbers in the O:S table (P3:S4)
# - Add a new "F-Measure" labeled block (rows 20-23) mirroring the existing
#   DT/RF/KNN/NB layout used by the ROC-AUC / Precision / Recall tables
# - Add two new charts that plot the new F-Measure data (a 3-D clustered
#   column chart off the O:S table, and a line chart off the new A21:E23 rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Update the existing "F-Measure" numeric columns (P3:S4) that already
#    existed in the O:S table but had stale values.
# ---------------------------------------------------------------------
$ws.Range("P3").Value2 = 87
$ws.Range("Q3").Value2 = 90
$ws.Range("R3").Value2 = 92
$ws.Range("S3").Value2 = 91

$ws.Range("P4").Value2 = 91
$ws.Range("Q4").Value2 = 92
$ws.Range("R4").Value2 = 94
$ws.Range("S4").Value2 = 94

# ---------------------------------------------------------------------
# 2. Add the new "F-Measure" block: header row 20, column headers row 21,
#    and the two data rows 22-23 (ROS / AVG ROS&ROS), mirroring the
#    existing A1:E4 "ROC-AUC" block layout + formatting.
# ---------------------------------------------------------------------

# Row 20: merged section title, formatted like A1:E1 / I1:M1 / O1:S1
$ws.Range("A1:E1").Copy() | Out-Null
$ws.Range("A20:E20").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A20").Value2 = "F-Measure"

# Row 21: blank corner + DT/RF/KNN/NB headers, formatted like row 2
$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A21:E21").PasteSpecial(-4122) | Out-Null
$ws.Range("B21").Value2 = "DT"
$ws.Range("C21").Value2 = "RF"
$ws.Range("D21").Value2 = "KNN"
$ws.Range("E21").Value2 = "NB"

# Row 22: "ROS" data row, formatted like row 3
$ws.Range("A3:E3").Copy() | Out-Null
$ws.Range("A22:E22").PasteSpecial(-4122) | Out-Null
$ws.Range("A22").Value2 = "ROS"
$ws.Range("B22").Value2 = 89
$ws.Range("C22").Value2 = 94
$ws.Range("D22").Value2 = 87
$ws.Range("E22").Value2 = 54

# Row 23: "AVG ROS&ROS" data row, formatted like row 4
$ws.Range("A4:E4").Copy() | Out-Null
$ws.Range("A23:E23").PasteSpecial(-4122) | Out-Null
$ws.Range("A23").Value2 = "AVG ROS&ROS"
$ws.Range("B23").Value2 = 92
$ws.Range("C23").Value2 = 94
$ws.Range("D23").Value2 = 89
$ws.Range("E23").Value2 = 56

$ws.Range("A20:E20").Merge() | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Move the active selection to reflect the new working area.
# ---------------------------------------------------------------------
$ws.Range("N24").Select() | Out-Null
try {
    $excel.ActiveWindow.ScrollRow = 10
} catch {
}

# ---------------------------------------------------------------------
# 4. Add the two new charts for the F-Measure data.
# ---------------------------------------------------------------------

# Chart 3: clustered column (3-D) chart built off the O:S "F-Measure" table,
# positioned to the right of the existing two charts.
$chart3Obj = $ws.ChartObjects().Add(560, 95, 320, 205)
$chart3 = $chart3Obj.Chart
$chart3.ChartType = 54   # xl3DColumnClustered
$chart3.SetSourceData($ws.Range("O2:S4"))
$chart3Obj.Name = "Chart 2"
$chart3Obj.Left = 560
$chart3Obj.Top = 95
$chart3Obj.Width = 320
$chart3Obj.Height = 205

# Chart 4: line chart with markers built off the new A21:E23 rows,
# positioned below the first two charts.
$chart4Obj = $ws.ChartObjects().Add(245, 290, 320, 205)
$chart4 = $chart4Obj.Chart
$chart4.ChartType = 65   # xlLineMarkers
$chart4.SetSourceData($ws.Range("A21:E23"))
$chart4Obj.Name = "Chart 6"
$chart4Obj.Left = 245
$chart4Obj.Top = 290
$chart4Obj.Width = 320
$chart4Obj.Height = 205
